$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.797.11"
$ws.Range("E2").Value = "  +3.24%  "

$ws.Range("D3").Value = "3.244.03"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'544.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "

$ws.Range("D6").Value = "'146.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("E10").Value = "  +2.64%  "

$ws.Range("E11").Value = "  -1.12%  "

$ws.Range("D12").Value = "3.804.59"
$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("E13").Value = "  -1.91%  "

$ws.Range("D14").Value = "'26.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").Value = "60.775.91"
$ws.Range("E16").Value = "  +3.10%  "

$ws.Range("D17").Value = "3.243.32"
$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("D18").Value = "'6.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").Value = "'13.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.25%  "

$ws.Range("D20").Value = "'8.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.98%  "

$ws.Range("D21").Value = "'377.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").Value = "'0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "'69.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").Value = "0.0₃0915"
$ws.Range("E28").Value = "  +5.62%  "

$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").Value = "'22.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").Value = "'6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.87%  "

$ws.Range("E32").Value = "  +4.09%  "

$ws.Range("E33").Value = "  +6.62%  "

$ws.Range("E34").Value = "  +4.65%  "

$ws.Range("D35").Value = "'158.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("E36").Value = "  +6.49%  "

$ws.Range("D37").Value = "'26.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.39%  "

$ws.Range("D38").Value = "2.803.22"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D39").Value = "'0.0721"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.92%  "

$ws.Range("E40").Value = "  +6.82%  "

$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("D43").Value = "'39.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "

$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").Value = "3.286.55"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("E46").Value = "  +3.06%  "

$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("D48").Value = "'21.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.17%  "

$ws.Range("D49").Value = "'6.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "'0.805"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.38%  "

$ws.Range("D51").Value = "'277.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.04%  "
